$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 721. This shifts the existing rows 721-797
# down to 722-798, matching the rest of the diff (which is just every
# subsequent row's content moving down by one row).
$ws.Rows.Item(721).Insert()

# Populate the newly inserted row 721 with the new price record.
$ws.Range("A721").Value = 6
$ws.Range("B721").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C721").Value = "Metropolitana"
$ws.Range("D721").Value = 45194
$ws.Range("E721").Value = 13
$ws.Range("F721").Value = 100112052
$ws.Range("G721").Value = "Albahaca"
$ws.Range("H721").Value = "Sin especificar"
$ws.Range("I721").Value = "Primera"
$ws.Range("J721").Value = 80
$ws.Range("K721").Value = 5000
$ws.Range("L721").Value = 5500
$ws.Range("M721").Value = 5156
$ws.Range("N721").Value = "`$/paquete"
$ws.Range("O721").Value = "Región de Arica y Parinacota"
$ws.Range("P721").Value = 5156
$ws.Range("Q721").Value = 1
$ws.Range("R721").Value = "Hortaliza"

# Match the date-number-format style used by the other date cells in column D.
$ws.Range("D721").NumberFormat = $ws.Range("D722").NumberFormat
